# Apply the updated "total_registros" survey counts and re-sort the
# empadronador list by that column in descending order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New (name, total_registros) pairs in the final, sorted order for rows 2-18.
# Rows 19-21 are unchanged.
$data = @(
    @("TOLENTINO VASQUEZ DIANA KATHERYN", 50),
    @("JULCA VALENZUELA CINTIA KARYN", 47),
    @("SANCHEZ CORTEZ LEYLA DIANA", 43),
    @("CARRILLO MARTÍNEZ HEIDY NAYELI", 41),
    @("VALER VEGA PATRICIA GERALDINE", 41),
    @("DE LA CRUZ BENITES RICHARD ALEXANDER", 40),
    @("YZQUIERDO CARHUATANTA LEYDY YANELA", 38),
    @("RODRIGUEZ RUBIO SANDRA MABEL", 38),
    @("PONCE VILLANUEVA CARMEN ISABEL", 36),
    @("ARENAS ZAVALA ANDYELA PATRICIA ISIDORA", 35),
    @("REYES RODRIGUEZ JEISSON STEVEN", 34),
    @("GASLAC GUTIERREZ FRANK JHORDY", 34),
    @("RUBIO MARIÑOS GISELA JUDITH", 32),
    @("PIERINA NAGIELLY SANDOVAL CONTRERAS", 31),
    @("CYNTHIA RODRIGUEZ LECCA", 28),
    @("SEGURA ASTO YAMILET ANTONELA", 25),
    @("RODRIGUEZ VASQUEZ WALTER", 20)
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $row++
}
